$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.668673753738403
$ws.Range("B1").Value = 1.929265975952148
$ws.Range("C1").Value = 1.751630187034607
$ws.Range("D1").Value = 2.503470659255981
$ws.Range("E1").Value = 3.077194213867188
